# Generate Report for Handoff
#
# Refresh the handoff timestamps / priority markers for the six source
# files that are currently "Ready for handoff" (rows 7, 8, 9, 11, 12, 14
# on each sheet). The "Overview" sheet tracks the latest handoff-xliff
# generation time per file (column G); the per-locale sheets ("zh-cn",
# "de-de") track the same handoff datetime (column H) and now also carry
# the "ht" (handoff type) priority marker (column E).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 12, 14)

# Overview sheet: bump the "Latest HO Xliff Generate Date" column (G).
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-09-03 18:32:32"
}

# zh-cn sheet: set Priority (E) and bump "Latest Handoff Datetime" (H).
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-03 18:32:27"
}

# de-de sheet: set Priority (E) and bump "Latest Handoff Datetime" (H).
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-03 18:32:32"
}
